# Add column for "Status as of July 4, 2025" on Sheet1 (with a dropdown list
# for the value), backed by a new hidden "DropdownOptions" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Create the hidden DropdownOptions sheet (placed right after Sheet1) ---
$ws1 = $wb.Worksheets.Item(1)
$dropSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$dropSheet.Name = "DropdownOptions"

$options = @("0% - 10%", "11% - 25%", "26% - 50%", "51% - 75%", "76% - 90%", "91% - 99%", "100%")
for ($i = 0; $i -lt $options.Length; $i++) {
    $cell = $dropSheet.Cells.Item($i + 1, 1)
    # Leading apostrophe forces text so values like "100%" aren't coerced
    # into a numeric percentage; resetting the style afterwards drops the
    # quote-prefix formatting so the cell keeps the plain default style.
    $cell.Value = "'" + $options[$i]
    $cell.Style = "Normal"
}

$dropSheet.Visible = $false

# --- 2. Add the new header + dropdown-enabled column on Sheet1 ---
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(1, 34).Value = "Status as of July 4, 2025"

# Clean up stray empty cells that previously held blank inline strings.
$ws.Cells.Item(2, 16).ClearContents()
$ws.Cells.Item(2, 31).ClearContents()
$ws.Cells.Item(2, 32).ClearContents()

$target = $ws.Cells.Item(2, 34)
$target.Validation.Add(3, 1, 1, '=DropdownOptions!$A$1:$A$7')
$target.Validation.ShowInput = $false
$target.Validation.ShowError = $false

# Keep Sheet1 as the active/selected sheet (matches pre-edit state).
$ws.Activate()

Write-Output "done"
